$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Slides"
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("Slides")

# Row 2 (S12)
$ws.Range("H2").Value = "Ajustes finos de padding/spacing para evitar corte; compactação de blocos e rodapé."
$ws.Range("J2").Value = "—"

# Row 13 (S18)
$ws.Range("F13").Value = "P1"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "2026-01-25"
$ws.Range("G13").ClearFormats()
$ws.Range("H13").Value = "Marcador de IC 95% redesenhado (mais limpo) + ajustes visuais menores."
$ws.Range("J13").Value = "—"

# Row 14 (S19)
$ws.Range("F14").Value = "P1"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "2026-01-25"
$ws.Range("G14").ClearFormats()
$ws.Range("H14").Value = "Tabela passou a usar estilo global table-medium (padding consistente)."
$ws.Range("J14").Value = "—"

# Row 16 (S23) - was only A:E populated (status DONE); now gets full row data
$ws.Range("E16").Value = "Atualizado"
$ws.Range("F16").Value = "P2"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "2026-01-25"
$ws.Range("G16").ClearFormats()
$ws.Range("H16").Value = "Redesign completo + conteúdo (TRIPOD/PROBAST, discriminação vs calibração, validação)."
$ws.Range("I16").Value = "TRIPOD/PROBAST"
$ws.Range("J16").Value = "—"

# Row 18 (S25)
$ws.Range("F18").Value = "P1"
$ws.Range("H18").Value = "Ajuste de paleta do card “Conceito” (header navy + chip gold) + ref PREVENT/PCE."
$ws.Range("J18").Value = "—"

# Row 21 (S51)
$ws.Range("E21").Value = "Atualizado"
$ws.Range("H21").Value = "Tabela ipsis litteris (Extremo) + chip “novo” e nota em rodapé."
$ws.Range("J21").Value = "—"

# Row 23 (S53)
$ws.Range("E23").Value = "Atualizado"
$ws.Range("F23").Value = "P1"
$ws.Range("H23").Value = "Atualização VESALIUS‑CV: inclusão de HR/IC (3‑point e 4‑point MACE)."
$ws.Range("J23").Value = "—"

# Row 25 (S55)
$ws.Range("E25").Value = "Atualizado"
$ws.Range("F25").Value = "P1"
$ws.Range("H25").Value = "Incluído lembrete “rate up” (Core GRADE) no contexto de viés de publicação."
$ws.Range("J25").Value = "—"

# Row 26 (S56)
$ws.Range("E26").Value = "Atualizado"
$ws.Range("H26").Value = "Linha AACE 2025 ajustada (meta LDL<70 + recomendações farmacológicas) + refs completas."
$ws.Range("J26").Value = "—"

# New row 30 (GRADE / S08)
$ws.Range("A30").Value = "GRADE"
$ws.Range("B30").Value = "S08"
$ws.Range("C30").Value = "Certeza da evidência (GRADE): domínio–resposta"
$ws.Range("D30").Value = "Fundamentos"
$ws.Range("E30").Value = "Atualizado"
$ws.Range("F30").Value = "P1"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "2026-01-25"
$ws.Range("G30").ClearFormats()
$ws.Range("H30").Value = "Ajuste de paleta (remoção de teal residual), alinhamento/spacing e padronização visual."
$ws.Range("I30").Value = "Core GRADE (visão geral)"
$ws.Range("J30").Value = "—"

# New row 31 (GRADE / S22)
$ws.Range("A31").Value = "GRADE"
$ws.Range("B31").Value = "S22"
$ws.Range("C31").Value = "Diretrizes: calculadoras de risco & papel do CAC"
$ws.Range("D31").Value = "Apêndice (PREVENT)"
$ws.Range("E31").Value = "Atualizado"
$ws.Range("F31").Value = "P2"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "2026-01-25"
$ws.Range("G31").ClearFormats()
$ws.Range("H31").Value = "Redesign completo (cards consistentes) + fontes e mensagem-chave (EtD)."
$ws.Range("I31").Value = "SBC 2025; ACC/AHA 2019; ESC 2021"
$ws.Range("J31").Value = "—"

# -----------------------------------------------------------------
# Sheet "Batches"
# -----------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Batches")

$wsB.Range("A5").NumberFormat = "@"
$wsB.Range("A5").Value = "2026-01-25"
$wsB.Range("A5").ClearFormats()
$wsB.Range("B5").Value = "Patch 2.7"
$wsB.Range("C5").Value = "P0/P1: navegação e ordem (main→metas→encerramento→apêndice); ajustes de paleta/padding; redesign PREVENT (S22–S23); atualização refs (AACE 2025, VESALIUS‑CV)."
$wsB.Range("D5").Value = "ZIP com arquivos modificados (HTML/CSS/JS + DASHBOARD + CHANGELOG)."
$wsB.Range("E5").Value = "Nav sem loop; h2 padronizado (2.85vw); table-medium; S51 ipsis litteris (Extremo)."
